# Add "hydrogen combined cycle" as a new power plant type, and rename the
# existing "hydrogen" entry to "hydrogen combustion turbine" on both the
# RQSD-BRQSD and RQSD-RQSD sheets.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("RQSD-BRQSD")
$ws3 = $wb.Worksheets.Item("RQSD-RQSD")

# --- RQSD-BRQSD sheet -------------------------------------------------
# Rename the existing "hydrogen" row (row 24) to "hydrogen combustion turbine"
$ws2.Range("A24").Value = "hydrogen combustion turbine"

# Give it the new look: black font color, vertically centered text.
$ws2.Range("A24").Font.Color = 0
$ws2.Range("A24").VerticalAlignment = -4108

# Add the new "hydrogen combined cycle" row right after it.
$ws2.Range("A25").Value = "hydrogen combined cycle"
$ws2.Range("B25").Value = 1

# Copy A24's formatting onto A25 (re-using the same cell style rather than
# re-deriving it from scratch for every cell).
$ws2.Range("A24").Copy()
$ws2.Range("A25").PasteSpecial(-4122)

# --- RQSD-RQSD sheet ----------------------------------------------------
$ws3.Range("A24").Value = "hydrogen combustion turbine"
$ws3.Range("A25").Value = "hydrogen combined cycle"
$ws3.Range("B25").Value = 1

# Re-use the already-formatted cell from the other sheet as the format
# source so we don't create a second, redundant cell-style entry.
$ws2.Range("A24").Copy()
$ws3.Range("A24").PasteSpecial(-4122)
$ws3.Range("A25").PasteSpecial(-4122)

$ws2.Range("B26").Select()
$ws3.Range("B26").Select()

# Restore the originally-active sheet/tab (the edits above only touched the
# per-sheet selection rectangle, not which tab is shown as active).
$wb.Worksheets.Item("About").Activate()
